# Auto update: 2025-12-06 02:00:05
# Re-labels row 4 / row 5 (Bitcoin <-> Coinbase swap caused by a shared-string
# reorder upstream) and refreshes the day's metrics (price/RSI/5d return/
# rule score/probability columns/final score) for every ticker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: now Coinbase Global, Inc. / COIN -------------------------------
$ws.Range("B4").Value = "Coinbase Global, Inc."
$ws.Range("C4").Value = "COIN"

# --- Row 5: now Bitcoin USD / BTC-USD --------------------------------------
$ws.Range("B5").Value = "Bitcoin USD"
$ws.Range("C5").Value = "BTC-USD"

# --- Row 2 (RIOT) updated metrics -------------------------------------------
$ws.Range("D2").Value = 14.93
$ws.Range("E2").Value = 57.2
$ws.Range("F2").Value = -7.42
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 53.7
$ws.Range("N2").Value = 49.16024380385575

# --- Row 3 (MARA) updated metrics -------------------------------------------
$ws.Range("D3").Value = 11.81
$ws.Range("E3").Value = 48.6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 48.9
$ws.Range("N3").Value = 49.16024380385575

# --- Row 4 (COIN) updated metrics -------------------------------------------
$ws.Range("D4").Value = 270.42
$ws.Range("E4").Value = 44.4
$ws.Range("F4").Value = -0.88
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 47.7
$ws.Range("N4").Value = 49.16024380385575

# --- Row 5 (BTC-USD) updated metrics ----------------------------------------
$ws.Range("D5").Value = 88706.96000000001
$ws.Range("E5").Value = 57
$ws.Range("F5").Value = -1.87
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 63
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 43
$ws.Range("K5").Value = 46.7
$ws.Range("N5").Value = 49.16024380385575

# --- Row 6 (MSTR) updated metrics -------------------------------------------
$ws.Range("D6").Value = 178.95
$ws.Range("E6").Value = 40.1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 36
$ws.Range("K6").Value = 42.7
$ws.Range("N6").Value = 49.16024380385575

Write-Output "edit applied"
